$d = $word.ActiveDocument

# 1. Replace the first paragraph's placeholder-ID text (which originally
#    spans two runs: the ID text run plus a trailing space run) with the
#    new ID text only. Matching the trailing space here lets the two runs
#    collapse into a single run with no residual whitespace, matching the
#    target markup exactly.
$null = $d.Content.Find.Execute("**ID__AFFARS_5339_topic_3__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5339_7201_90__ID**", 2)

# 2. Update the first paragraph's paragraph formatting: increase the left
#    indent from 6pt (120 twips) to 11.25pt (225 twips) and add a
#    paragraph border with 5pt spacing on all four sides (no line drawn).
$p = $d.Paragraphs(1)
$pf = $p.Range.ParagraphFormat
$pf.LeftIndent = 11.25

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
